# Apply the "updated classes and docs" change:
#  - add a new ScheduleName / MyAstroSchedule row to the Configuration sheet
#  - resize the Configuration columns to fit the new content
#  - make Configuration the active/selected sheet (was Entries)
#  - move the selection on Entries to F9, and keep Configuration's selection on B7

$wb = $excel.ActiveWorkbook

$wsEntries = $wb.Worksheets.Item("Entries")
$wsConfig  = $wb.Worksheets.Item("Configuration")

# --- Configuration sheet: add the new ScheduleName row ---
$wsConfig.Range("A7").Value = "ScheduleName"
$wsConfig.Range("B7").Value = "MyAstroSchedule"

# Widen the columns so the new, longer values fit (mirrors the bestFit resize
# Excel performs automatically when the column content grows).
$wsConfig.Columns.Item(1).ColumnWidth = 13.65
$wsConfig.Columns.Item(2).ColumnWidth = 15.85

# --- Selection / active sheet bookkeeping ---
$wsEntries.Range("F9").Select()

$wsConfig.Activate()
$wsConfig.Range("B7").Select()
